$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Target cluster changes from ECs to MuSCs, with refreshed TPM stats)
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3971766666666667
$ws.Range("H2").Value = 1.19153
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.443866
$ws.Range("N2").Value = 0.887732
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1762932183266667
$ws.Range("R2").Value = 1.05775930996
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove the now-obsolete rows 3 and 4 (old ECs->MuSCs / ECs->Resolving-Mac pairs)
$ws.Rows("3:4").Delete()
